# indicadores.xlsx - update "ÁREA TECH" certification/enrolment figures and
# refresh the active sheet/selection to match where the author left off
# working (streamlit cache-reload button + consultant table cleanup).

$wb = $excel.ActiveWorkbook

$wsGestion = $wb.Worksheets.Item("ÁREA GESTIÓN CORPORATIVA")
$wsTech    = $wb.Worksheets.Item("ÁREA TECH")

# --- ÁREA TECH: "Master Profesional en IA" / "Certificación SAP" block ---
# The right-hand mini table used to re-use the "50 ALUMNOS" caption; it now
# gets its own "21 ALUMNOS" label, and the enrolment/closed counters below
# it are reset/updated.
$wsTech.Range("F25").Value = "21 ALUMNOS"
$wsTech.Range("E27").Value = 0

$wsTech.Range("F28").Value = 1
$wsTech.Range("F29").Value = 0
$wsTech.Range("F30").Value = 0
$wsTech.Range("F33").Value = 0
$wsTech.Range("F34").Value = 0.89
$wsTech.Range("F35").ClearContents()
$wsTech.Range("F36").Value = 26

# --- ÁREA TECH: certifications table (consultores) cleanup ---
$wsTech.Range("F39").Value = 1
$wsTech.Range("F40").Value = 10
$wsTech.Range("F41").Value = 0
$wsTech.Range("F42").Value = 0
$wsTech.Range("F43").Value = 0

# --- View state: user ended up with "ÁREA TECH" as the active/selected
# sheet (scrolled down to the certifications block), while
# "ÁREA GESTIÓN CORPORATIVA" lost its former selection.
$wsGestion.Range("E29").Select() | Out-Null
$wsTech.Select() | Out-Null
$wsTech.Range("K37").Select() | Out-Null
